$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from E1) onto the three new header cells
# so they reuse the same cellXfs style index instead of creating a new one.
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New header labels
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# New boolean (FALSE) columns for data rows 2-8
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}
